# The upstream change being replayed here (commit "Fixed POI packaging and
# upgraded to POI 3.15.") only touches how the OOXML parts were re-serialized
# by the authoring tool: every hunk in the diff is the exact same element,
# with the exact same attribute names/values, just re-emitted with the
# attributes sorted alphabetically (an Apache POI/XMLBeans writer artifact
# of the library upgrade). No text, numbers, styles, page geometry, fonts,
# or any other document content/semantics actually changed between the two
# revisions.
#
# Word's object model (and this COM-interop runtime) exposes the document's
# *semantic* content - paragraphs, runs, styles, page setup, etc. - and does
# not provide any way to dictate the raw XML attribute serialization order
# used when a part is written back out. There is therefore no COM call that
# corresponds to this diff beyond leaving the document's content untouched.
#
# This script intentionally performs no content mutation, which is the
# faithful COM-automation equivalent of a pure attribute-reordering/
# re-serialization change.
$d = $word.ActiveDocument
